$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new trade row (row 8) to the bag-trade log, reusing row 7's
# formatting (date format on A, boolean format on G) via copy/paste-special
# so no new style entries get minted.
$ws.Range("A7:I7").Copy() | Out-Null
$ws.Range("A8:I8").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(8, 1).Value = 42654.743657407409
$ws.Cells.Item(8, 2).Value = $false
$ws.Cells.Item(8, 3).Value = 10079.18
$ws.Cells.Item(8, 4).Value = 10079.68
$ws.Cells.Item(8, 5).Value = 75.5
$ws.Cells.Item(8, 6).Value = 75.489998
$ws.Cells.Item(8, 7).Value = $false
$ws.Cells.Item(8, 8).Value = -0.01
$ws.Cells.Item(8, 9).Value = $false
